# ---------------------------------------------------------------------------
# Fills in the marksheet with the graded quiz results (previously the sheet
# held placeholder "Absent" data) and drops the now-unused 3rd question
# block (columns G/H).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# --- Summary block (rows 10-12) -------------------------------------------

# Give the row-label cells in column A the same "absoluteStyle" formatting
# that row 9's labels already use, then fill in the labels.
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)

$ws.Range("A10").Value = "No."
$ws.Range("A11").Value = "Marking"
$ws.Range("A12").Value = "Total"

# Right / Wrong / Not-Attempt / Max counts for the graded attempt.
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 28

# Marking scheme per question.
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

# Totals.
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "76/112"

# --- Per-question "Student Ans" columns (A for block 1, D for block 1b) ---
# Column A rows 16-40 mirror block 1's student answers; a blank means the
# question was not attempted (kept at the existing "normalStyle" s=7).
# Where the student's answer matches the "Correct Ans" column (B) the cell
# is painted with the "correctStyle" (copy format from B10, s=5); otherwise
# with the "incorrectStyle" (copy format from C10, s=6).

$correct = @{
    16 = "Option A"; 17 = "Option D"; 18 = "Option B"; 19 = "Option C";
    21 = "Option C"; 22 = "Option D"; 23 = "Option D";
    25 = "Option A"; 26 = "Option C"; 28 = "Option D";
    32 = "Option C"; 33 = "Option D"; 34 = "Option B";
    36 = "Option A"; 37 = "Option A"; 38 = "Option A"; 39 = "Option D"; 40 = "Option D"
}
$incorrect = @{
    27 = "Option D"; 31 = "Option C"; 35 = "Option B"
}

$ws.Range("B10").Copy()
foreach ($row in $correct.Keys) {
    $ws.Range("A$row").PasteSpecial(-4122)
}
$ws.Range("C10").Copy()
foreach ($row in $incorrect.Keys) {
    $ws.Range("A$row").PasteSpecial(-4122)
}

foreach ($row in $correct.Keys) {
    $ws.Range("A$row").Value = $correct[$row]
}
foreach ($row in $incorrect.Keys) {
    $ws.Range("A$row").Value = $incorrect[$row]
}

# Column D rows 16-18 (block 1b student answers); row 17 was answered wrong.
$ws.Range("B10").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("D17").PasteSpecial(-4122)

$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option D"
$ws.Range("D18").Value = "Option D"

# --- Drop the now-unused 3rd question block (rows 15-21, columns G/H), and
# the rest of the 1b block's student/correct-answer pairs (rows 19-40,
# columns D/E) which no longer apply now that the block only has 3 questions.

$ws.Range("D19:E40").Clear()
$ws.Range("G:H").EntireColumn.Delete()
